$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$full = $tr.Characters(1, $tr.Length)
$full.Text = "Two-Column Layout"
